# "maj template comment à la fin" — move the "Comment" column (data-info
# sheet columns J:L, 5 rows: header / French label / type / format / example)
# to the end, so the column order becomes:
#   before: J=Comment, K=ArrayType, L=Result
#   after:  J=ArrayType, K=Result,  L=Comment
# i.e. ArrayType and Result each shift one column left, and Comment's whole
# column of content moves to column L.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: this COM host's `Range.Value` getter has a reflection quirk that
# returns the property descriptor text instead of invoking the getter —
# `Value2` (and `Text`) work correctly for both read and write, so use that.
for ($row = 1; $row -le 5; $row++) {
    $commentCell   = $ws.Range("J$row")
    $arrayTypeCell = $ws.Range("K$row")
    $resultCell    = $ws.Range("L$row")

    $commentValue   = $commentCell.Value2
    $arrayTypeValue = $arrayTypeCell.Value2
    $resultValue    = $resultCell.Value2

    $commentCell.Value2   = $arrayTypeValue   # J: ArrayType
    $arrayTypeCell.Value2 = $resultValue      # K: Result
    $resultCell.Value2    = $commentValue     # L: Comment
}
